$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 44444
$ws.Range("J3").Value = 44444
$ws.Range("L3").Value = 44444
$ws.Range("N3").Value = -44672
$ws.Range("H18").Value = 516.55554
$ws.Range("I18").Value = 516.55554
$ws.Range("K18").Value = 516.55554
$ws.Range("M18").Value = -232.55554
$ws.Range("H19").Value = 562.4
$ws.Range("I19").Value = 379.66666
$ws.Range("J19").Value = 640.7143
$ws.Range("K19").Value = 379.66666
$ws.Range("L19").Value = 640.7143
$ws.Range("M19").Value = -204.66666
$ws.Range("N19").Value = -990.7143
$ws.Range("H40").Value = 1294
$ws.Range("I40").Value = 1242.5
$ws.Range("J40").Value = 1500
$ws.Range("K40").Value = 1242.5
$ws.Range("L40").Value = 1500
$ws.Range("M40").Value = -1067.5
$ws.Range("N40").Value = -1850
$ws.Range("H62").Value = 4639.6
$ws.Range("I62").Value = 4099.5
$ws.Range("K62").Value = 4099.5
$ws.Range("M62").Value = -3475.5
$ws.Range("H65").Value = 4639.6
$ws.Range("I65").Value = 4099.5
$ws.Range("K65").Value = 20497.5
$ws.Range("M65").Value = -17377.5
$ws.Range("H69").Value = 5001
$ws.Range("I69").Value = 5001
$ws.Range("K69").Value = 15003
$ws.Range("M69").Value = -14129
$ws.Range("H70").Value = 2294.6667
$ws.Range("I70").Value = 1828
$ws.Range("J70").Value = 2994.6667
$ws.Range("K70").Value = 5484
$ws.Range("L70").Value = 8984.000100000001
$ws.Range("M70").Value = -5214
$ws.Range("N70").Value = -9524.000100000001
$ws.Range("H72").Value = 5001
$ws.Range("I72").Value = 5001
$ws.Range("K72").Value = 45009
$ws.Range("M72").Value = -40641
$ws.Range("H73").Value = 2294.6667
$ws.Range("I73").Value = 1828
$ws.Range("J73").Value = 2994.6667
$ws.Range("K73").Value = 5484
$ws.Range("L73").Value = 8984.000100000001
$ws.Range("M73").Value = -4548
$ws.Range("N73").Value = -10856.0001
$ws.Range("H76").Value = 13927
$ws.Range("I76").Value = 10699.667
$ws.Range("K76").Value = 10699.667
$ws.Range("M76").Value = -10384.667
$ws.Range("H79").Value = 13927
$ws.Range("I79").Value = 10699.667
$ws.Range("K79").Value = 10699.667
$ws.Range("M79").Value = -9607.666999999999
$ws.Range("H93").Value = 19999.5
$ws.Range("J93").Value = 19999.5
$ws.Range("L93").Value = 19999.5
$ws.Range("N93").Value = -24991.5
$ws.Range("H98").Value = 2777.647
$ws.Range("I98").Value = 2419.6072
$ws.Range("K98").Value = 2419.6072
$ws.Range("M98").Value = -921.6071999999999
$ws.Range("H99").Value = 1106.5714
$ws.Range("I99").Value = 925.1
$ws.Range("J99").Value = 1560.25
$ws.Range("K99").Value = 2775.3
$ws.Range("L99").Value = 4680.75
$ws.Range("M99").Value = -1277.3
$ws.Range("N99").Value = -7676.75
$ws.Range("H102").Value = 44444
$ws.Range("J102").Value = 44444
$ws.Range("L102").Value = 44444
$ws.Range("N102").Value = -50934
$ws.Range("H106").Value = 5699.4
$ws.Range("I106").Value = 5699.4
$ws.Range("K106").Value = 5699.4
$ws.Range("M106").Value = -5068.4
$ws.Range("H111").Value = 1756.25
$ws.Range("I111").Value = 1165
$ws.Range("K111").Value = 3495
$ws.Range("M111").Value = -428
$ws.Range("H116").Value = 2781.4
$ws.Range("J116").Value = 3335.3635
$ws.Range("L116").Value = 3335.3635
$ws.Range("N116").Value = -10219.3635
$ws.Range("H122").Value = 2777.647
$ws.Range("I122").Value = 2419.6072
$ws.Range("K122").Value = 7258.821599999999
$ws.Range("M122").Value = -4808.821599999999
$ws.Range("H132").Value = 13981.174
$ws.Range("I132").Value = 10517.286
$ws.Range("K132").Value = 31551.858
$ws.Range("M132").Value = -29021.858
$ws.Range("H137").Value = 9928.48
$ws.Range("I137").Value = 1440.3334
$ws.Range("K137").Value = 4321.0002
$ws.Range("M137").Value = -1771.0002
$ws.Range("H138").Value = 3550.4312
$ws.Range("I138").Value = 1057.6
$ws.Range("K138").Value = 3172.8
$ws.Range("M138").Value = 1967.2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1267.48
$ws.Range("I2").Value = 1207.5294
$ws.Range("J2").Value = 1394.875
$ws.Range("K2").Value = 1207.5294
$ws.Range("L2").Value = 1394.875
$ws.Range("M2").Value = -1094.5294
$ws.Range("N2").Value = -1620.875
$ws.Range("H55").Value = 48659.332
$ws.Range("J55").Value = 48659.332
$ws.Range("L55").Value = 48659.332
$ws.Range("N55").Value = -49289.332
$ws.Range("H74").Value = 18545.732
$ws.Range("I74").Value = 1546.3889
$ws.Range("J74").Value = 44044.75
$ws.Range("K74").Value = 1546.3889
$ws.Range("L74").Value = 44044.75
$ws.Range("M74").Value = -672.3888999999999
$ws.Range("N74").Value = -45792.75
$ws.Range("H77").Value = 18545.732
$ws.Range("I77").Value = 1546.3889
$ws.Range("J77").Value = 44044.75
$ws.Range("K77").Value = 7731.9445
$ws.Range("L77").Value = 220223.75
$ws.Range("M77").Value = -3363.9445
$ws.Range("N77").Value = -228959.75
$ws.Range("H88").Value = 1997
$ws.Range("I88").Value = 1893
$ws.Range("K88").Value = 1893
$ws.Range("M88").Value = -1487
$ws.Range("H91").Value = 1997
$ws.Range("I91").Value = 1893
$ws.Range("K91").Value = 1893
$ws.Range("M91").Value = -489
$ws.Range("H94").Value = 30000
$ws.Range("J94").Value = 30000
$ws.Range("L94").Value = 30000
$ws.Range("N94").Value = -31802
$ws.Range("H106").Value = 49999.5
$ws.Range("I106").Value = 15000
$ws.Range("J106").Value = 84999
$ws.Range("K106").Value = 15000
$ws.Range("L106").Value = 84999
$ws.Range("M106").Value = -13738
$ws.Range("N106").Value = -87523
$ws.Range("H116").Value = 1267.48
$ws.Range("I116").Value = 1207.5294
$ws.Range("J116").Value = 1394.875
$ws.Range("K116").Value = 1207.5294
$ws.Range("L116").Value = 1394.875
$ws.Range("M116").Value = 1086.4706
$ws.Range("N116").Value = -5982.875
$ws.Range("H130").Value = 30750
$ws.Range("J130").Value = 30750
$ws.Range("L130").Value = 30750
$ws.Range("N130").Value = -40790
$ws.Range("H132").Value = 12527222
$ws.Range("I132").Value = 3550
$ws.Range("K132").Value = 10650
$ws.Range("M132").Value = -8120
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1267.48
$ws.Range("I3").Value = 1207.5294
$ws.Range("J3").Value = 1394.875
$ws.Range("K3").Value = 1207.5294
$ws.Range("L3").Value = 1394.875
$ws.Range("M3").Value = -1093.5294
$ws.Range("N3").Value = -1622.875
$ws.Range("H80").Value = 1147
$ws.Range("I80").Value = 358.75
$ws.Range("J80").Value = 1389.5385
$ws.Range("K80").Value = 358.75
$ws.Range("L80").Value = 1389.5385
$ws.Range("M80").Value = 639.25
$ws.Range("N80").Value = -3385.5385
$ws.Range("H83").Value = 1147
$ws.Range("I83").Value = 358.75
$ws.Range("J83").Value = 1389.5385
$ws.Range("K83").Value = 1793.75
$ws.Range("L83").Value = 6947.692500000001
$ws.Range("M83").Value = 3198.25
$ws.Range("N83").Value = -16931.6925
$ws.Range("H103").Value = 22599.6
$ws.Range("J103").Value = 22599.6
$ws.Range("L103").Value = 22599.6
$ws.Range("N103").Value = -24943.6
$ws.Range("H105").Value = 1135.05
$ws.Range("I105").Value = 1182.5385
$ws.Range("J105").Value = 1046.8572
$ws.Range("K105").Value = 1182.5385
$ws.Range("L105").Value = 1046.8572
$ws.Range("M105").Value = 564.4614999999999
$ws.Range("N105").Value = -4540.8572
$ws.Range("H107").Value = 1987.2858
$ws.Range("I107").Value = 1970.3334
$ws.Range("K107").Value = 1970.3334
$ws.Range("M107").Value = -50.33339999999998
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 3399.6667
$ws.Range("I2").Value = 3333
$ws.Range("J2").Value = 3466.3333
$ws.Range("K2").Value = 3333
$ws.Range("L2").Value = 3466.3333
$ws.Range("M2").Value = -3220
$ws.Range("N2").Value = -3692.3333
$ws.Range("H3").Value = 3226
$ws.Range("I3").Value = 3226
$ws.Range("K3").Value = 3226
$ws.Range("M3").Value = -3113
$ws.Range("H4").Value = 400
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 400
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 400
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -624
$ws.Range("H5").Value = 1118.1818
$ws.Range("I5").Value = 307.85715
$ws.Range("J5").Value = 2536.25
$ws.Range("K5").Value = 307.85715
$ws.Range("L5").Value = 2536.25
$ws.Range("M5").Value = -195.85715
$ws.Range("N5").Value = -2760.25
$ws.Range("H8").Value = 2701.25
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 2701.25
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 2701.25
$ws.Range("M8").Value = -2981.25
$ws.Range("H10").Value = 101005.5
$ws.Range("J10").Value = 201860.8
$ws.Range("L10").Value = 201860.8
$ws.Range("N10").Value = -202138.8
$ws.Range("H12").Value = 455250.28
$ws.Range("I12").Value = 907.25
$ws.Range("K12").Value = 907.25
$ws.Range("M12").Value = -737.25
$ws.Range("H13").Value = 1163.3334
$ws.Range("J13").Value = 1727.5
$ws.Range("L13").Value = 1727.5
$ws.Range("N13").Value = -2005.5
$ws.Range("H14").Value = 2336.6667
$ws.Range("I14").Value = 510
$ws.Range("K14").Value = 510
$ws.Range("M14").Value = -340
$ws.Range("H15").Value = 397.5
$ws.Range("J15").Value = 397.5
$ws.Range("L15").Value = 397.5
$ws.Range("N15").Value = -737.5
$ws.Range("H22").Value = 2591.2307
$ws.Range("I22").Value = 1743.2222
$ws.Range("K22").Value = 1743.2222
$ws.Range("M22").Value = -1393.2222
$ws.Range("H31").Value = 17490.264
$ws.Range("I31").Value = 1911.6
$ws.Range("J31").Value = 75910.25
$ws.Range("K31").Value = 1911.6
$ws.Range("L31").Value = 75910.25
$ws.Range("M31").Value = -1616.6
$ws.Range("N31").Value = -76500.25
$ws.Range("H34").Value = 17490.264
$ws.Range("I34").Value = 1911.6
$ws.Range("J34").Value = 75910.25
$ws.Range("K34").Value = 1911.6
$ws.Range("L34").Value = 75910.25
$ws.Range("M34").Value = -1709.6
$ws.Range("N34").Value = -76314.25
$ws.Range("H43").Value = 19500
$ws.Range("J43").Value = 19500
$ws.Range("L43").Value = 19500
$ws.Range("N43").Value = -19868
$ws.Range("H58").Value = 19516.61
$ws.Range("I58").Value = 11583.75
$ws.Range("J58").Value = 23747.467
$ws.Range("K58").Value = 11583.75
$ws.Range("L58").Value = 23747.467
$ws.Range("M58").Value = -11380.75
$ws.Range("N58").Value = -24153.467
$ws.Range("H62").Value = 7749.6665
$ws.Range("I62").Value = 8266.333000000001
$ws.Range("K62").Value = 8266.333000000001
$ws.Range("M62").Value = -7642.333000000001
$ws.Range("H65").Value = 7749.6665
$ws.Range("I65").Value = 8266.333000000001
$ws.Range("K65").Value = 41331.665
$ws.Range("M65").Value = -38211.665
$ws.Range("H74").Value = 37999
$ws.Range("J74").Value = 37999
$ws.Range("L74").Value = 37999
$ws.Range("N74").Value = -39747
$ws.Range("H77").Value = 37999
$ws.Range("J77").Value = 37999
$ws.Range("L77").Value = 113997
$ws.Range("N77").Value = -122733
$ws.Range("H99").Value = 4649.8335
$ws.Range("I99").Value = 5224.75
$ws.Range("K99").Value = 5224.75
$ws.Range("M99").Value = -3726.75
$ws.Range("H101").Value = 19500
$ws.Range("J101").Value = 19500
$ws.Range("L101").Value = 19500
$ws.Range("N101").Value = -25990
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H126").Value = 4649.8335
$ws.Range("I126").Value = 5224.75
$ws.Range("K126").Value = 15674.25
$ws.Range("M126").Value = -13204.25
$ws.Range("H136").Value = 19516.61
$ws.Range("I136").Value = 11583.75
$ws.Range("J136").Value = 23747.467
$ws.Range("K136").Value = 34751.25
$ws.Range("L136").Value = 71242.401
$ws.Range("M136").Value = -32201.25
$ws.Range("N136").Value = -76342.401
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 1178.4
$ws.Range("I75").Value = 797.3333
$ws.Range("K75").Value = 2391.9999
$ws.Range("M75").Value = -1393.9999
$ws.Range("H76").Value = 2500
$ws.Range("I76").Value = 2500
$ws.Range("K76").Value = 7500
$ws.Range("M76").Value = -7117
$ws.Range("H78").Value = 1178.4
$ws.Range("I78").Value = 797.3333
$ws.Range("K78").Value = 7175.9997
$ws.Range("M78").Value = -2183.9997
$ws.Range("H79").Value = 2500
$ws.Range("I79").Value = 2500
$ws.Range("K79").Value = 7500
$ws.Range("M79").Value = -6174
$ws.Range("H86").Value = 693
$ws.Range("I86").Value = 686.2727
$ws.Range("J86").Value = 711.5
$ws.Range("K86").Value = 2058.8181
$ws.Range("L86").Value = 2134.5
$ws.Range("M86").Value = -872.8181
$ws.Range("N86").Value = -4506.5
$ws.Range("H89").Value = 693
$ws.Range("I89").Value = 686.2727
$ws.Range("J89").Value = 711.5
$ws.Range("K89").Value = 6176.454299999999
$ws.Range("L89").Value = 6403.5
$ws.Range("M89").Value = -248.4542999999994
$ws.Range("N89").Value = -18259.5
$ws.Range("H94").Value = 998.5
$ws.Range("I94").Value = 998.5
$ws.Range("K94").Value = 2995.5
$ws.Range("M94").Value = -2319.5
$ws.Range("H97").Value = 450
$ws.Range("J97").Value = 500
$ws.Range("L97").Value = 1500
$ws.Range("N97").Value = -2492
$ws.Range("H122").Value = 17938920
$ws.Range("I122").Value = 31145712
$ws.Range("J122").Value = 4732129
$ws.Range("K122").Value = 280311408
$ws.Range("L122").Value = 42589161
$ws.Range("M122").Value = -280308958
$ws.Range("N122").Value = -42594061
$ws.Range("H131").Value = 1466.13
$ws.Range("J131").Value = 1485.0532
$ws.Range("L131").Value = 4455.1596
$ws.Range("N131").Value = -14535.1596
$ws.Range("H137").Value = 2789.1428
$ws.Range("I137").Value = 1591.1666
$ws.Range("J137").Value = 9977
$ws.Range("K137").Value = 4773.4998
$ws.Range("L137").Value = 29931
$ws.Range("M137").Value = 326.5002000000004
$ws.Range("N137").Value = -40131
$ws.Range("H140").Value = 2273.3333
$ws.Range("I140").Value = 2710
$ws.Range("K140").Value = 8130
$ws.Range("M140").Value = -2950
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 96.25
$ws.Range("I2").Value = 100.86667
$ws.Range("J2").Value = 88.55556
$ws.Range("K2").Value = 100.86667
$ws.Range("L2").Value = 88.55556
$ws.Range("M2").Value = 12.13333
$ws.Range("N2").Value = -314.55556
$ws.Range("H26").Value = 22989.75
$ws.Range("J26").Value = 22989.75
$ws.Range("L26").Value = 22989.75
$ws.Range("N26").Value = -23549.75
$ws.Range("H50").Value = 22989.75
$ws.Range("J50").Value = 22989.75
$ws.Range("L50").Value = 22989.75
$ws.Range("N50").Value = -23985.75
$ws.Range("H55").Value = 19665
$ws.Range("J55").Value = 19665
$ws.Range("L55").Value = 19665
$ws.Range("N55").Value = -20319
$ws.Range("H62").Value = 19500
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 19500
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H93").Value = 34995
$ws.Range("J93").Value = 34995
$ws.Range("L93").Value = 34995
$ws.Range("N93").Value = -38739
$ws.Range("H96").Value = 10000
$ws.Range("J96").Value = 10000
$ws.Range("L96").Value = 10000
$ws.Range("N96").Value = -15492
$ws.Range("H102").Value = 6874.1304
$ws.Range("I102").Value = 7504.4736
$ws.Range("J102").Value = 3880
$ws.Range("K102").Value = 7504.4736
$ws.Range("L102").Value = 3880
$ws.Range("M102").Value = -5882.4736
$ws.Range("N102").Value = -7124
$ws.Range("H105").Value = 30069
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H113").Value = 2170.6667
$ws.Range("I113").Value = 2183.2222
$ws.Range("J113").Value = 2133
$ws.Range("K113").Value = 2183.2222
$ws.Range("L113").Value = 2133
$ws.Range("M113").Value = -13.22220000000016
$ws.Range("N113").Value = -6473
$ws.Range("H123").Value = 52387
$ws.Range("J123").Value = 52387
$ws.Range("L123").Value = 52387
$ws.Range("N123").Value = -57287
$ws.Range("H126").Value = 7269
$ws.Range("I126").Value = 7345.364
$ws.Range("K126").Value = 22036.092
$ws.Range("M126").Value = -19566.092
$ws.Range("H132").Value = 13639.833
$ws.Range("I132").Value = 4000
$ws.Range("K132").Value = 12000
$ws.Range("M132").Value = -9470
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("H16").Value = 2454.6924
$ws.Range("I16").Value = 2259.6365
$ws.Range("J16").Value = 3527.5
$ws.Range("K16").Value = 2259.6365
$ws.Range("L16").Value = 3527.5
$ws.Range("M16").Value = -2089.6365
$ws.Range("N16").Value = -3867.5
$ws.Range("H19").Value = 3559.1667
$ws.Range("I19").Value = 713.75
$ws.Range("K19").Value = 713.75
$ws.Range("M19").Value = -543.75
$ws.Range("H34").Value = 30000
$ws.Range("I34").Value = 30000
$ws.Range("K34").Value = 30000
$ws.Range("M34").Value = -29828
$ws.Range("H40").Value = 627.7143
$ws.Range("I40").Value = 565.6667
$ws.Range("J40").Value = 1000
$ws.Range("K40").Value = 565.6667
$ws.Range("L40").Value = 1000
$ws.Range("M40").Value = -429.6667
$ws.Range("N40").Value = -1272
$ws.Range("H61").Value = 1389.1333
$ws.Range("I61").Value = 1467.3572
$ws.Range("J61").Value = 294
$ws.Range("K61").Value = 1467.3572
$ws.Range("L61").Value = 294
$ws.Range("M61").Value = -1265.3572
$ws.Range("N61").Value = -698
$ws.Range("J63").Value = 20333.334
$ws.Range("L63").Value = 20333.334
$ws.Range("N63").Value = -21831.334
$ws.Range("J66").Value = 20333.334
$ws.Range("L66").Value = 61000.00199999999
$ws.Range("N66").Value = -68488.00199999999
$ws.Range("H68").Value = 6399.8
$ws.Range("I68").Value = 6399.8
$ws.Range("K68").Value = 6399.8
$ws.Range("M68").Value = -5650.8
$ws.Range("H71").Value = 6399.8
$ws.Range("I71").Value = 6399.8
$ws.Range("K71").Value = 31999
$ws.Range("M71").Value = -28255
$ws.Range("H104").Value = 20513.334
$ws.Range("J104").Value = 20513.334
$ws.Range("L104").Value = 20513.334
$ws.Range("N104").Value = -27501.334
$ws.Range("H106").Value = 13046
$ws.Range("J106").Value = 13046
$ws.Range("L106").Value = 13046
$ws.Range("N106").Value = -15570
$ws.Range("H113").Value = 1389.1333
$ws.Range("I113").Value = 1467.3572
$ws.Range("J113").Value = 294
$ws.Range("K113").Value = 1467.3572
$ws.Range("L113").Value = 294
$ws.Range("M113").Value = 702.6428000000001
$ws.Range("N113").Value = -4634
$ws.Range("H136").Value = 258186.11
$ws.Range("I136").Value = 35945.168
$ws.Range("K136").Value = 107835.504
$ws.Range("M136").Value = -105285.504
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 5512750
$ws.Range("I3").Value = 11000000
$ws.Range("J3").Value = 25500.5
$ws.Range("K3").Value = 11000000
$ws.Range("L3").Value = 25500.5
$ws.Range("M3").Value = -10999886
$ws.Range("N3").Value = -25728.5
$ws.Range("H34").Value = 217000
$ws.Range("J34").Value = 225499.5
$ws.Range("L34").Value = 225499.5
$ws.Range("N34").Value = -225905.5
$ws.Range("H62").Value = 13073.869
$ws.Range("J62").Value = 13086.533
$ws.Range("L62").Value = 13086.533
$ws.Range("N62").Value = -14334.533
$ws.Range("H65").Value = 13073.869
$ws.Range("J65").Value = 13086.533
$ws.Range("L65").Value = 65432.66499999999
$ws.Range("N65").Value = -71672.66499999999
$ws.Range("H107").Value = 1027.2
$ws.Range("I107").Value = 1359.8889
$ws.Range("J107").Value = 528.1667
$ws.Range("K107").Value = 4079.6667
$ws.Range("L107").Value = 1584.5001
$ws.Range("M107").Value = -2159.6667
$ws.Range("N107").Value = -5424.5001
$ws.Range("H126").Value = 18371.455
$ws.Range("I126").Value = 19174.857
$ws.Range("K126").Value = 57524.571
$ws.Range("M126").Value = -55054.571
$ws.Range("H135").Value = 93583.336
$ws.Range("J135").Value = 93583.336
$ws.Range("L135").Value = 93583.336
$ws.Range("N135").Value = -103723.336
$ws.Range("H136").Value = 14745.218
$ws.Range("I136").Value = 1956.25
$ws.Range("K136").Value = 5868.75
$ws.Range("M136").Value = -3318.75
